$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3076923076923077
$ws.Range("C2").Value = 0.3846153846153846
$ws.Range("P2").Value = 0.2692307692307692
$ws.Range("S2").Value = 0.03846153846153846

# Row 3
$ws.Range("J3").Value = 0.1
$ws.Range("P3").Value = 0.9

# Row 4
$ws.Range("J4").Value = 0.3333333333333333
$ws.Range("P4").Value = 0.6666666666666666

# Row 6
$ws.Range("J6").Value = 0.09090909090909091
$ws.Range("Q6").Value = 0.1818181818181818
$ws.Range("R6").Value = 0.1818181818181818
$ws.Range("S6").Value = 0.5454545454545454

# Row 7
$ws.Range("B7").Value = 0.375
$ws.Range("Q7").Value = 0.125
$ws.Range("R7").Value = 0.25
$ws.Range("S7").Value = 0.25

# Row 8
$ws.Range("B8").Value = 0.05
$ws.Range("F8").Value = 0.1
$ws.Range("J8").Value = 0.15
$ws.Range("Q8").Value = 0.3
$ws.Range("R8").Value = 0.15
$ws.Range("S8").Value = 0.25

# Row 9
$ws.Range("B9").Value = 0.2142857142857143
$ws.Range("F9").Value = 0.07142857142857142
$ws.Range("J9").Value = 0.1428571428571428
$ws.Range("Q9").Value = 0.2142857142857143
$ws.Range("R9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.2142857142857143

# Row 10
$ws.Range("B10").Value = 0.1052631578947368
$ws.Range("D10").Value = 0.02631578947368421
$ws.Range("F10").Value = 0.05263157894736842
$ws.Range("J10").Value = 0.2105263157894737
$ws.Range("O10").Value = 0.008771929824561403
$ws.Range("Q10").Value = 0.2982456140350877
$ws.Range("R10").Value = 0.08771929824561403
$ws.Range("S10").Value = 0.2105263157894737

# Row 11
$ws.Range("G11").Value = 0.25
$ws.Range("K11").Value = 0.25
$ws.Range("L11").Value = 0.5

# Row 12
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.3333333333333333

# Row 13
$ws.Range("G13").Value = 0.5
$ws.Range("J13").Value = 0.5

# Row 15
$ws.Range("H15").Value = 0.1578947368421053
$ws.Range("I15").Value = 0.1578947368421053
$ws.Range("J15").Value = 0.5263157894736842
$ws.Range("K15").Value = 0.05263157894736842
$ws.Range("O15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.05263157894736842

# Row 16
$ws.Range("H16").Value = 0.1052631578947368
$ws.Range("J16").Value = 0.631578947368421
$ws.Range("K16").Value = 0.1578947368421053
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.05263157894736842

# Row 17
$ws.Range("H17").Value = 0.1333333333333333
$ws.Range("I17").Value = 0.1111111111111111
$ws.Range("J17").Value = 0.6
$ws.Range("K17").Value = 0.06666666666666667
$ws.Range("O17").Value = 0.04444444444444445
$ws.Range("S17").Value = 0.04444444444444445

# Row 18
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("O18").Value = 0.1111111111111111

# Row 19
$ws.Range("H19").Value = 0.1304347826086956
$ws.Range("I19").Value = 0.08695652173913043
$ws.Range("J19").Value = 0.4782608695652174
$ws.Range("M19").Value = 0.04347826086956522
$ws.Range("O19").Value = 0.2173913043478261
$ws.Range("S19").Value = 0.04347826086956522
